$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("H2").Value = 5.57
$ws1.Range("L2").Value = 1.1

$ws1.Range("H3").Value = 4.38
$ws1.Range("L3").Value = 1

$ws1.Range("H4").Value = 3.6
$ws1.Range("L4").Value = 0.92

$ws1.Range("D5").Value = 63
$ws1.Range("H5").Value = 2.39
$ws1.Range("L5").Value = 1.15

$ws1.Range("D6").Value = 66
$ws1.Range("H6").Value = 1.34
$ws1.Range("L6").Value = 0.83

$ws1.Range("H7").Value = 0.36
$ws1.Range("I7").Value = "High"
$ws1.Range("J7").Value = "Urgent"
$ws1.Range("L7").Value = 1

$ws1.Range("H8").Value = 0
$ws1.Range("L8").Value = 1.15

$ws1.Range("L9").Value = 0.92

$ws1.Range("L10").Value = 1.17

$ws1.Range("L11").Value = 1.15

$ws1.Range("L12").Value = 0.86

$ws1.Range("L13").Value = 0.95

$ws1.Range("L14").Value = 0.89

$ws1.Range("L15").Value = 0.87

$ws1.Range("L16").Value = 0.86

$ws1.Range("L17").Value = 1.18

# --- Sheet: Summary ---
# These cells hold numbers that are stored as text in the original workbook,
# so force a text number format before assigning to keep them as text cells.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "1028"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "508"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "244"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "58"
